$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1970
$ws.Range("I28").Value = 1644.1428
$ws.Range("K28").Value = 1644.1428
$ws.Range("M28").Value = -1159.1428
$ws.Range("H51").Value = 9594.619000000001
$ws.Range("J51").Value = 9946.684999999999
$ws.Range("L51").Value = 9946.684999999999
$ws.Range("N51").Value = -10914.685
$ws.Range("H98").Value = 2301.4
$ws.Range("I98").Value = 1875.25
$ws.Range("K98").Value = 1875.25
$ws.Range("M98").Value = -377.25
$ws.Range("H121").Value = 4747.8696
$ws.Range("J121").Value = 4747.8696
$ws.Range("L121").Value = 14243.6088
$ws.Range("N121").Value = -17737.6088
$ws.Range("H122").Value = 2301.4
$ws.Range("I122").Value = 1875.25
$ws.Range("K122").Value = 5625.75
$ws.Range("M122").Value = -3175.75
$ws.Range("H132").Value = 7178.857
$ws.Range("I132").Value = 1579.1724
$ws.Range("J132").Value = 34244
$ws.Range("K132").Value = 4737.5172
$ws.Range("L132").Value = 102732
$ws.Range("M132").Value = -2207.5172
$ws.Range("N132").Value = -107792
$ws.Range("H137").Value = 30306464
$ws.Range("I137").Value = 2586.8333
$ws.Range("K137").Value = 7760.499899999999
$ws.Range("M137").Value = -5210.499899999999
$ws.Range("H141").Value = 8976.76
$ws.Range("I141").Value = 7431.1
$ws.Range("K141").Value = 22293.3
$ws.Range("M141").Value = -17113.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5322.8276
$ws.Range("I61").Value = 4849.7896
$ws.Range("K61").Value = 4849.7896
$ws.Range("M61").Value = -4637.7896
$ws.Range("H74").Value = 45455868
$ws.Range("I74").Value = 71429496
$ws.Range("J74").Value = 2024.75
$ws.Range("K74").Value = 71429496
$ws.Range("L74").Value = 2024.75
$ws.Range("M74").Value = -71428622
$ws.Range("N74").Value = -3772.75
$ws.Range("H77").Value = 45455868
$ws.Range("I77").Value = 71429496
$ws.Range("J77").Value = 2024.75
$ws.Range("K77").Value = 357147480
$ws.Range("L77").Value = 10123.75
$ws.Range("M77").Value = -357143112
$ws.Range("N77").Value = -18859.75
$ws.Range("H101").Value = 28750
$ws.Range("J101").Value = 28750
$ws.Range("L101").Value = 28750
$ws.Range("N101").Value = -35240
$ws.Range("H136").Value = 5322.8276
$ws.Range("I136").Value = 4849.7896
$ws.Range("K136").Value = 14549.3688
$ws.Range("M136").Value = -11999.3688

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 15109320
$ws.Range("I94").Value = 19026284
$ws.Range("J94").Value = 1028.1428
$ws.Range("K94").Value = 19026284
$ws.Range("L94").Value = 1028.1428
$ws.Range("M94").Value = -19025833
$ws.Range("N94").Value = -1930.1428
$ws.Range("H99").Value = 1097971.6
$ws.Range("I99").Value = 1158802.9
$ws.Range("K99").Value = 1158802.9
$ws.Range("M99").Value = -1157304.9
$ws.Range("H105").Value = 5999.3335
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1270.3
$ws.Range("I16").Value = 963.6
$ws.Range("J16").Value = 1577
$ws.Range("K16").Value = 963.6
$ws.Range("L16").Value = 1577
$ws.Range("M16").Value = -676.6
$ws.Range("N16").Value = -2151
$ws.Range("H31").Value = 12501954
$ws.Range("I31").Value = 14926982
$ws.Range("K31").Value = 14926982
$ws.Range("M31").Value = -14926687
$ws.Range("H34").Value = 12501954
$ws.Range("I34").Value = 14926982
$ws.Range("K34").Value = 14926982
$ws.Range("M34").Value = -14926780
$ws.Range("H58").Value = 3539.1428
$ws.Range("I58").Value = 2793.3333
$ws.Range("K58").Value = 2793.3333
$ws.Range("M58").Value = -2590.3333
$ws.Range("H113").Value = 1270.3
$ws.Range("I113").Value = 963.6
$ws.Range("J113").Value = 1577
$ws.Range("K113").Value = 963.6
$ws.Range("L113").Value = 1577
$ws.Range("M113").Value = 1206.4
$ws.Range("N113").Value = -5917
$ws.Range("H134").Value = 2632.0908
$ws.Range("I134").Value = 2376.9412
$ws.Range("J134").Value = 3499.6
$ws.Range("K134").Value = 7130.823600000001
$ws.Range("L134").Value = 10498.8
$ws.Range("M134").Value = -4595.823600000001
$ws.Range("N134").Value = -15568.8
$ws.Range("H136").Value = 3539.1428
$ws.Range("I136").Value = 2793.3333
$ws.Range("K136").Value = 8379.999899999999
$ws.Range("M136").Value = -5829.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 8249.75
$ws.Range("J74").Value = 11499.5
$ws.Range("L74").Value = 34498.5
$ws.Range("N74").Value = -36620.5
$ws.Range("H77").Value = 8249.75
$ws.Range("J77").Value = 11499.5
$ws.Range("L77").Value = 103495.5
$ws.Range("N77").Value = -114103.5
$ws.Range("H132").Value = 1201
$ws.Range("I132").Value = 989.35297
$ws.Range("K132").Value = 8904.176730000001
$ws.Range("M132").Value = -6374.176730000001
$ws.Range("H134").Value = 5228
$ws.Range("I134").Value = 1689.1428
$ws.Range("K134").Value = 5067.428400000001
$ws.Range("M134").Value = 2.571599999999307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3326199.8
$ws.Range("I80").Value = 5541916.5
$ws.Range("J80").Value = 2625
$ws.Range("K80").Value = 5541916.5
$ws.Range("L80").Value = 2625
$ws.Range("M80").Value = -5540918.5
$ws.Range("N80").Value = -4621
$ws.Range("H83").Value = 3326199.8
$ws.Range("I83").Value = 5541916.5
$ws.Range("J83").Value = 2625
$ws.Range("K83").Value = 27709582.5
$ws.Range("L83").Value = 13125
$ws.Range("M83").Value = -27704590.5
$ws.Range("N83").Value = -23109
$ws.Range("H113").Value = 2359529.2
$ws.Range("I113").Value = 2696447.8
$ws.Range("K113").Value = 2696447.8
$ws.Range("M113").Value = -2694277.8
$ws.Range("H140").Value = 86777.336
$ws.Range("J140").Value = 86777.336
$ws.Range("L140").Value = 86777.336
$ws.Range("N140").Value = -97137.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 15000000
$ws.Range("I33").Value = 15000000
$ws.Range("K33").Value = 15000000
$ws.Range("M33").Value = -14999710
$ws.Range("H46").Value = 4100.6875
$ws.Range("J46").Value = 4921.5557
$ws.Range("L46").Value = 4921.5557
$ws.Range("N46").Value = -5297.5557
$ws.Range("H136").Value = 4034.1836
$ws.Range("I136").Value = 2819.9167
$ws.Range("J136").Value = 7396.769
$ws.Range("K136").Value = 8459.750100000001
$ws.Range("L136").Value = 22190.307
$ws.Range("M136").Value = -5909.750100000001
$ws.Range("N136").Value = -27290.307

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2348.762
$ws.Range("I132").Value = 2351.889
$ws.Range("J132").Value = 2330
$ws.Range("K132").Value = 7055.667
$ws.Range("L132").Value = 6990
$ws.Range("M132").Value = -4525.667
$ws.Range("N132").Value = -12050
$ws.Range("H136").Value = 4009.5356
$ws.Range("I136").Value = 2707.348
$ws.Range("J136").Value = 9999.6
$ws.Range("K136").Value = 8122.044
$ws.Range("L136").Value = 29998.8
$ws.Range("M136").Value = -5572.044
$ws.Range("N136").Value = -35098.8
